$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the affected cells keep their original text (inline string) data type
# even when the new value looks numeric (e.g. "1.004"), matching how the sheet
# originally stored these columns as text.
$cells = @(
    @{ Ref = "D2"; Val = "27.621.03" },
    @{ Ref = "E2"; Val = "  -0.22%  " },
    @{ Ref = "D3"; Val = "1.845.18" },
    @{ Ref = "E3"; Val = "  -0.20%  " },
    @{ Ref = "D4"; Val = "1.004" },
    @{ Ref = "E4"; Val = "  +0.16%  " },
    @{ Ref = "D5"; Val = "315.44" },
    @{ Ref = "E5"; Val = "  +0.87%  " },
    @{ Ref = "D6"; Val = "1.003" },
    @{ Ref = "E6"; Val = "  +0.18%  " },
    @{ Ref = "D7"; Val = "0.4307" },
    @{ Ref = "E7"; Val = "  +0.86%  " },
    @{ Ref = "D8"; Val = "0.3688" },
    @{ Ref = "E8"; Val = "  +1.55%  " },
    @{ Ref = "B9"; Val = "OKB" },
    @{ Ref = "C9"; Val = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb" },
    @{ Ref = "D9"; Val = "45.16" },
    @{ Ref = "E9"; Val = "  +0.76%  " },
    @{ Ref = "B10"; Val = "Dogecoin" },
    @{ Ref = "C10"; Val = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge" },
    @{ Ref = "D10"; Val = "0.07313" },
    @{ Ref = "E10"; Val = "  -0.11%  " },
    @{ Ref = "B11"; Val = "Polygon" },
    @{ Ref = "C11"; Val = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic" },
    @{ Ref = "D11"; Val = "0.8753" },
    @{ Ref = "E11"; Val = "  -0.10%  " },
    @{ Ref = "B12"; Val = "Solana" },
    @{ Ref = "C12"; Val = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol" },
    @{ Ref = "D12"; Val = "20.98" },
    @{ Ref = "E12"; Val = "  +1.70%  " },
    @{ Ref = "B13"; Val = "WrappedEther" },
    @{ Ref = "C13"; Val = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth" },
    @{ Ref = "D13"; Val = "1.877.63" },
    @{ Ref = "E13"; Val = "  -3.08%  " },
    @{ Ref = "B14"; Val = "Polkadot" },
    @{ Ref = "C14"; Val = "https://coinranking.com/coin/25W7FG7om+polkadot-dot" },
    @{ Ref = "D14"; Val = "5.471" },
    @{ Ref = "E14"; Val = "  +2.81%  " },
    @{ Ref = "B15"; Val = "Chainlink" },
    @{ Ref = "C15"; Val = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link" },
    @{ Ref = "D15"; Val = "6.602" },
    @{ Ref = "E15"; Val = "  +1.21%  " },
    @{ Ref = "B16"; Val = "TRON" },
    @{ Ref = "C16"; Val = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx" },
    @{ Ref = "D16"; Val = "0.06964" },
    @{ Ref = "E16"; Val = "  +0.80%  " },
    @{ Ref = "B17"; Val = "BinanceUSD" },
    @{ Ref = "C17"; Val = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd" },
    @{ Ref = "D17"; Val = "1.006" },
    @{ Ref = "E17"; Val = "  +0.24%  " },
    @{ Ref = "B18"; Val = "Litecoin" },
    @{ Ref = "C18"; Val = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc" },
    @{ Ref = "D18"; Val = "81.36" },
    @{ Ref = "E18"; Val = "  +1.84%  " },
    @{ Ref = "B19"; Val = "ShibaInu" },
    @{ Ref = "C19"; Val = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib" },
    @{ Ref = "D19"; Val = "0.000009059" },
    @{ Ref = "E19"; Val = "  +0.29%  " },
    @{ Ref = "B20"; Val = "Dai" },
    @{ Ref = "C20"; Val = "https://coinranking.com/coin/MoTuySvg7+dai-dai" },
    @{ Ref = "D20"; Val = "1.004" },
    @{ Ref = "E20"; Val = "  -0.01%  " },
    @{ Ref = "B21"; Val = "Avalanche" },
    @{ Ref = "C21"; Val = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax" },
    @{ Ref = "D21"; Val = "15.59" },
    @{ Ref = "E21"; Val = "  +1.67%  " },
    @{ Ref = "B22"; Val = "WrappedBTC" },
    @{ Ref = "C22"; Val = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc" },
    @{ Ref = "D22"; Val = "27.688.31" },
    @{ Ref = "E22"; Val = "  -0.05%  " },
    @{ Ref = "B23"; Val = "Uniswap" },
    @{ Ref = "C23"; Val = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni" },
    @{ Ref = "D23"; Val = "5.084" },
    @{ Ref = "E23"; Val = "  +2.40%  " },
    @{ Ref = "B24"; Val = "Cosmos" },
    @{ Ref = "C24"; Val = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom" },
    @{ Ref = "D24"; Val = "10.98" },
    @{ Ref = "E24"; Val = "  +5.80%  " },
    @{ Ref = "B25"; Val = "WrappedliquidstakedEther2.0" },
    @{ Ref = "C25"; Val = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth" },
    @{ Ref = "D25"; Val = "2.075.95" },
    @{ Ref = "E25"; Val = "  -3.34%  " },
    @{ Ref = "B26"; Val = "Toncoin" },
    @{ Ref = "C26"; Val = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton" },
    @{ Ref = "D26"; Val = "1.982" },
    @{ Ref = "E26"; Val = "  +0.83%  " },
    @{ Ref = "B27"; Val = "Monero" },
    @{ Ref = "C27"; Val = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr" },
    @{ Ref = "D27"; Val = "154.23" },
    @{ Ref = "E27"; Val = "  -0.06%  " },
    @{ Ref = "B28"; Val = "EthereumClassic" },
    @{ Ref = "C28"; Val = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc" },
    @{ Ref = "D28"; Val = "18.98" },
    @{ Ref = "E28"; Val = "  +0.92%  " },
    @{ Ref = "B29"; Val = "InternetComputer(DFINITY)" },
    @{ Ref = "C29"; Val = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp" },
    @{ Ref = "D29"; Val = "5.322" },
    @{ Ref = "E29"; Val = "  +1.04%  " },
    @{ Ref = "B30"; Val = "BitcoinCash" },
    @{ Ref = "C30"; Val = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch" },
    @{ Ref = "D30"; Val = "115.85" },
    @{ Ref = "E30"; Val = "  -4.69%  " },
    @{ Ref = "B31"; Val = "LidoDAOToken" },
    @{ Ref = "C31"; Val = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo" },
    @{ Ref = "D31"; Val = "1.883" },
    @{ Ref = "E31"; Val = "  +0.89%  " },
    @{ Ref = "B32"; Val = "Stellar" },
    @{ Ref = "C32"; Val = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm" },
    @{ Ref = "D32"; Val = "0.08898" },
    @{ Ref = "E32"; Val = "  -0.29%  " },
    @{ Ref = "B33"; Val = "ImmutableX" },
    @{ Ref = "C33"; Val = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx" },
    @{ Ref = "D33"; Val = "0.7842" },
    @{ Ref = "E33"; Val = "  +3.05%  " },
    @{ Ref = "B34"; Val = "Filecoin" },
    @{ Ref = "C34"; Val = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil" },
    @{ Ref = "D34"; Val = "4.609" },
    @{ Ref = "E34"; Val = "  +1.79%  " },
    @{ Ref = "B35"; Val = "HuobiToken" },
    @{ Ref = "C35"; Val = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht" },
    @{ Ref = "D35"; Val = "2.989" },
    @{ Ref = "E35"; Val = "  +0.76%  " },
    @{ Ref = "B36"; Val = "ARBITRUM" },
    @{ Ref = "C36"; Val = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb" },
    @{ Ref = "D36"; Val = "1.166" },
    @{ Ref = "E36"; Val = "  +5.84%  " },
    @{ Ref = "B37"; Val = "Frax" },
    @{ Ref = "C37"; Val = "https://coinranking.com/coin/KfWtaeV1W+frax-frax" },
    @{ Ref = "D37"; Val = "1.003" },
    @{ Ref = "E37"; Val = "  +0.24%  " },
    @{ Ref = "D38"; Val = "0.05437" },
    @{ Ref = "E38"; Val = "  +0.90%  " },
    @{ Ref = "B39"; Val = "TrustWalletToken" },
    @{ Ref = "C39"; Val = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt" },
    @{ Ref = "D39"; Val = "1.108" },
    @{ Ref = "E39"; Val = "  +1.38%  " },
    @{ Ref = "B40"; Val = "VeChain" },
    @{ Ref = "C40"; Val = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet" },
    @{ Ref = "D40"; Val = "0.01962" },
    @{ Ref = "E40"; Val = "  +1.54%  " },
    @{ Ref = "B41"; Val = "MXToken" },
    @{ Ref = "C41"; Val = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx" },
    @{ Ref = "D41"; Val = "2.835" },
    @{ Ref = "E41"; Val = "  +0.36%  " },
    @{ Ref = "B42"; Val = "TheSandbox" },
    @{ Ref = "C42"; Val = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand" },
    @{ Ref = "D42"; Val = "0.5178" },
    @{ Ref = "E42"; Val = "  +1.93%  " },
    @{ Ref = "B43"; Val = "Algorand" },
    @{ Ref = "C43"; Val = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo" },
    @{ Ref = "D43"; Val = "0.1691" },
    @{ Ref = "E43"; Val = "  +2.28%  " },
    @{ Ref = "B44"; Val = "FraxShare" },
    @{ Ref = "C44"; Val = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs" },
    @{ Ref = "D44"; Val = "6.757" },
    @{ Ref = "E44"; Val = "  -0.31%  " },
    @{ Ref = "B45"; Val = "Aptos" },
    @{ Ref = "C45"; Val = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt" },
    @{ Ref = "D45"; Val = "8.645" },
    @{ Ref = "E45"; Val = "  +3.51%  " },
    @{ Ref = "B46"; Val = "EnergySwap" },
    @{ Ref = "C46"; Val = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" },
    @{ Ref = "D46"; Val = "10.61" },
    @{ Ref = "E46"; Val = "  +2.90%  " },
    @{ Ref = "D47"; Val = "0.4794" },
    @{ Ref = "E47"; Val = "  +2.42%  " },
    @{ Ref = "B48"; Val = "Cronos" },
    @{ Ref = "C48"; Val = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro" },
    @{ Ref = "D48"; Val = "0.06545" },
    @{ Ref = "E48"; Val = "  +0.00%  " },
    @{ Ref = "B49"; Val = "Quant" },
    @{ Ref = "C49"; Val = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt" },
    @{ Ref = "D49"; Val = "106.32" },
    @{ Ref = "E49"; Val = "  +1.16%  " },
    @{ Ref = "D50"; Val = "1.003" },
    @{ Ref = "E50"; Val = "  +0.21%  " },
    @{ Ref = "E51"; Val = "  +2.71%  " }
)

foreach ($item in $cells) {
    $cell = $ws.Range($item.Ref)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Val
    $cell.Style = "Normal"
}
